$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.304.70"
$ws.Range("E2").Value = "  -0.17%  "

# Row 3
$ws.Range("D3").Value = "1.910.97"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'0.721"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.18%  "

# Row 6
$ws.Range("D6").Value = "'254.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.64%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").Value = "'40.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.73%  "

# Row 9
$ws.Range("D9").Value = "'0.369"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.36%  "

# Row 10
$ws.Range("D10").Value = "'52.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.01%  "

# Row 11
$ws.Range("E11").Value = "  +6.88%  "

# Row 12
$ws.Range("E12").Value = "  -0.77%  "

# Row 13
$ws.Range("D13").Value = "2.186.82"
$ws.Range("E13").Value = "  -0.03%  "

# Row 14
$ws.Range("D14").Value = "'12.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.39%  "

# Row 15
$ws.Range("D15").Value = "'0.723"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.19%  "

# Row 16
$ws.Range("E16").Value = "  +1.53%  "

# Row 17
$ws.Range("D17").Value = "1.907.48"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18
$ws.Range("D18").Value = "35.289.66"
$ws.Range("E18").Value = "  -0.13%  "

# Row 19
$ws.Range("E19").Value = "  +2.48%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0855"
$ws.Range("E20").Value = "  +3.12%  "

# Row 21
$ws.Range("D21").Value = "'243.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.79%  "

# Row 23
$ws.Range("D23").Value = "'5.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.57%  "

# Row 24
$ws.Range("E24").Value = "  +0.21%  "

# Row 25
$ws.Range("E25").Value = "  +4.36%  "

# Row 26
$ws.Range("E26").Value = "  +3.90%  "

# Row 27
$ws.Range("D27").Value = "'167.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "

# Row 28
$ws.Range("D28").Value = "'8.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.42%  "

# Row 29
$ws.Range("D29").Value = "'18.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.86%  "

# Row 30
$ws.Range("E30").Value = "  +4.80%  "

# Row 31
$ws.Range("D31").Value = "4.129.70"
$ws.Range("E31").Value = "  +19.48%  "

# Row 32
$ws.Range("E32").Value = "  +5.03%  "

# Row 33
$ws.Range("E33").Value = "  +14.46%  "

# Row 34
$ws.Range("D34").Value = "'1.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +23.28%  "

# Row 35
$ws.Range("D35").Value = "'0.0587"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.70%  "

# Row 36
$ws.Range("D36").Value = "'4.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.82%  "

# Row 37
$ws.Range("E37").Value = "  -0.78%  "

# Row 38
$ws.Range("D38").Value = "'0.907"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.35%  "

# Row 39
$ws.Range("E39").Value = "  -0.35%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0218"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.89%  "

# Row 41
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'17.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.97%  "

# Row 42
$ws.Range("D42").Value = "'96.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.71%  "

# Row 43
$ws.Range("E43").Value = "  +0.69%  "

# Row 44
$ws.Range("E44").Value = "  +1.47%  "

# Row 45
$ws.Range("D45").Value = "1.336.14"
$ws.Range("E45").Value = "  -0.42%  "

# Row 46
$ws.Range("D46").Value = "'2.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.38%  "

# Row 47
$ws.Range("E47").Value = "  +1.00%  "

# Row 48
$ws.Range("D48").Value = "'6.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.87%  "

# Row 49
$ws.Range("E49").Value = "  -0.81%  "

# Row 50
$ws.Range("D50").Value = "'45.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.57%  "

# Row 51
$ws.Range("D51").Value = "'11.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.91%  "
